$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1996.4286
$ws.Range("J17").Value = 1996.4286
$ws.Range("L17").Value = 5989.2858
$ws.Range("N17").Value = -6325.2858
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H86").Value = 4039.6
$ws.Range("I86").Value = 4265.6665
$ws.Range("J86").Value = 3983.0833
$ws.Range("K86").Value = 4265.6665
$ws.Range("L86").Value = 3983.0833
$ws.Range("M86").Value = -3142.6665
$ws.Range("N86").Value = -6229.0833
$ws.Range("H89").Value = 4039.6
$ws.Range("I89").Value = 4265.6665
$ws.Range("J89").Value = 3983.0833
$ws.Range("K89").Value = 21328.3325
$ws.Range("L89").Value = 19915.4165
$ws.Range("M89").Value = -15712.3325
$ws.Range("N89").Value = -31147.4165
$ws.Range("H106").Value = 27530.385
$ws.Range("I106").Value = 28741.334
$ws.Range("K106").Value = 28741.334
$ws.Range("M106").Value = -28110.334
$ws.Range("H131").Value = 5380.231
$ws.Range("I131").Value = 689.4
$ws.Range("K131").Value = 2068.2
$ws.Range("M131").Value = 2971.8
$ws.Range("H132").Value = 2661.2222
$ws.Range("I132").Value = 2661.2222
$ws.Range("K132").Value = 7983.6666
$ws.Range("M132").Value = -5453.6666
$ws.Range("H137").Value = 2471.6
$ws.Range("I137").Value = 1260
$ws.Range("K137").Value = 3780
$ws.Range("M137").Value = -1230
$ws.Range("H141").Value = 3483.375
$ws.Range("I141").Value = 2266.7144
$ws.Range("K141").Value = 6800.1432
$ws.Range("M141").Value = -1620.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4998.3335
$ws.Range("I2").Value = 4997.5
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 4997.5
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -4884.5
$ws.Range("N2").Value = -5226
$ws.Range("H32").Value = 12719.909
$ws.Range("I32").Value = 8115.1875
$ws.Range("K32").Value = 8115.1875
$ws.Range("M32").Value = -7828.1875
$ws.Range("H74").Value = 1264.4568
$ws.Range("I74").Value = 982.8933
$ws.Range("K74").Value = 982.8933
$ws.Range("M74").Value = -108.8933
$ws.Range("H77").Value = 1264.4568
$ws.Range("I77").Value = 982.8933
$ws.Range("K77").Value = 4914.4665
$ws.Range("M77").Value = -546.4664999999995
$ws.Range("H116").Value = 4998.3335
$ws.Range("I116").Value = 4997.5
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4997.5
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -2703.5
$ws.Range("N116").Value = -9588
$ws.Range("H122").Value = 2774.75
$ws.Range("I122").Value = 2449.6667
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 7349.000100000001
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -4899.000100000001
$ws.Range("N122").Value = -16150
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4998.3335
$ws.Range("I3").Value = 4997.5
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 4997.5
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -4883.5
$ws.Range("N3").Value = -5228
$ws.Range("H105").Value = 4437.1577
$ws.Range("I105").Value = 3017.0557
$ws.Range("K105").Value = 3017.0557
$ws.Range("M105").Value = -1270.0557
$ws.Range("H134").Value = 2495.4211
$ws.Range("I134").Value = 2161
$ws.Range("K134").Value = 6483
$ws.Range("M134").Value = -3948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 7777
$ws.Range("J21").Value = 7777
$ws.Range("L21").Value = 7777
$ws.Range("N21").Value = -8247
$ws.Range("H31").Value = 5103.524
$ws.Range("I31").Value = 2330
$ws.Range("K31").Value = 2330
$ws.Range("M31").Value = -2035
$ws.Range("H34").Value = 5103.524
$ws.Range("I34").Value = 2330
$ws.Range("K34").Value = 2330
$ws.Range("M34").Value = -2128
$ws.Range("H58").Value = 3906.0908
$ws.Range("I58").Value = 1887.125
$ws.Range("J58").Value = 5059.7856
$ws.Range("K58").Value = 1887.125
$ws.Range("L58").Value = 5059.7856
$ws.Range("M58").Value = -1684.125
$ws.Range("N58").Value = -5465.7856
$ws.Range("H99").Value = 11732.469
$ws.Range("J99").Value = 12915.737
$ws.Range("L99").Value = 12915.737
$ws.Range("N99").Value = -15911.737
$ws.Range("H122").Value = 2592.261
$ws.Range("I122").Value = 2502.5
$ws.Range("J122").Value = 2915.4
$ws.Range("K122").Value = 7507.5
$ws.Range("L122").Value = 8746.200000000001
$ws.Range("M122").Value = -5057.5
$ws.Range("N122").Value = -13646.2
$ws.Range("H126").Value = 11732.469
$ws.Range("J126").Value = 12915.737
$ws.Range("L126").Value = 38747.211
$ws.Range("N126").Value = -43687.211
$ws.Range("H136").Value = 3906.0908
$ws.Range("I136").Value = 1887.125
$ws.Range("J136").Value = 5059.7856
$ws.Range("K136").Value = 5661.375
$ws.Range("L136").Value = 15179.3568
$ws.Range("M136").Value = -3111.375
$ws.Range("N136").Value = -20279.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 503.08334
$ws.Range("I11").Value = 148.14285
$ws.Range("K11").Value = 444.42855
$ws.Range("M11").Value = -304.42855
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H39").Value = 1750
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 1666.6666
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 4999.9998
$ws.Range("M39").Value = -5706
$ws.Range("N39").Value = -5587.9998
$ws.Range("H49").Value = 326.44446
$ws.Range("I49").Value = 326.44446
$ws.Range("K49").Value = 979.33338
$ws.Range("M49").Value = -823.33338
$ws.Range("H60").Value = 259.5
$ws.Range("I60").Value = 259.5
$ws.Range("K60").Value = 778.5
$ws.Range("M60").Value = -527.5
$ws.Range("H75").Value = 360
$ws.Range("J75").Value = 325
$ws.Range("L75").Value = 975
$ws.Range("N75").Value = -2971
$ws.Range("H78").Value = 360
$ws.Range("J78").Value = 325
$ws.Range("L78").Value = 2925
$ws.Range("N78").Value = -12909

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5001.6665
$ws.Range("I80").Value = 5001.6665
$ws.Range("K80").Value = 5001.6665
$ws.Range("M80").Value = -4003.6665
$ws.Range("H83").Value = 5001.6665
$ws.Range("I83").Value = 5001.6665
$ws.Range("K83").Value = 25008.3325
$ws.Range("M83").Value = -20016.3325
$ws.Range("H122").Value = 95320.45
$ws.Range("I122").Value = 3346.4
$ws.Range("J122").Value = 171965.5
$ws.Range("K122").Value = 10039.2
$ws.Range("L122").Value = 515896.5
$ws.Range("M122").Value = -7589.200000000001
$ws.Range("N122").Value = -520796.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1225.7693
$ws.Range("I16").Value = 1529.3
$ws.Range("J16").Value = 214
$ws.Range("K16").Value = 1529.3
$ws.Range("L16").Value = 214
$ws.Range("M16").Value = -1359.3
$ws.Range("N16").Value = -554
$ws.Range("H55").Value = 205.1579
$ws.Range("J55").Value = 289.33334
$ws.Range("L55").Value = 289.33334
$ws.Range("N55").Value = -635.33334
$ws.Range("H93").Value = 1048.0476
$ws.Range("I93").Value = 1048.0476
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1048.0476
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 199.9523999999999
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2050.7778
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2050.7778
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 1977.5
$ws.Range("I100").Value = 1722.3334
$ws.Range("K100").Value = 3444.6668
$ws.Range("M100").Value = -2903.6668
$ws.Range("H126").Value = 2386.1
$ws.Range("I126").Value = 1693.7142
$ws.Range("K126").Value = 5081.142599999999
$ws.Range("M126").Value = -2611.142599999999
